$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 10.30359681186211
$ws.Range("D2").Value = 4.307170786526501
$ws.Range("E2").Value = 13.75688765628082
$ws.Range("F2").Value = 20.81465131686366
$ws.Range("G2").Value = 22.81094877805262
$ws.Range("H2").Value = 12.77604502816191
$ws.Range("L2").Value = 9.942087910354596
$ws.Range("M2").Value = 59.00858740866833
$ws.Range("O2").Value = 18.65194293546658

$ws.Range("C3").Value = 10.52589910774381
$ws.Range("D3").Value = 4.282789355988567
$ws.Range("E3").Value = 13.57444447906548
$ws.Range("F3").Value = 21.05762261584053
$ws.Range("G3").Value = 23.21076216917748
$ws.Range("H3").Value = 12.9247932795092
$ws.Range("L3").Value = 9.781369041244401
$ws.Range("M3").Value = 55.53906150769269
$ws.Range("O3").Value = 18.92538273464109

$ws.Range("C4").Value = 10.66938621641625
$ws.Range("D4").Value = 4.268324493625647
$ws.Range("E4").Value = 13.46955324654705
$ws.Range("F4").Value = 21.22088918171962
$ws.Range("G4").Value = 23.48032654777419
$ws.Range("H4").Value = 13.02087133413101
$ws.Range("L4").Value = 9.685213236214285
$ws.Range("M4").Value = 53.28623388914487
$ws.Range("O4").Value = 19.10366012448817

$ws.Range("C5").Value = 10.72959278598603
$ws.Range("D5").Value = 4.262561409186408
$ws.Range("E5").Value = 13.42863555048814
$ws.Range("F5").Value = 21.29088390885052
$ws.Range("G5").Value = 23.59602688376929
$ws.Range("H5").Value = 13.06121086225098
$ws.Range("L5").Value = 9.64670958811625
$ws.Range("M5").Value = 52.33756258268932
$ws.Range("O5").Value = 19.17888308190894

$ws.Range("C6").Value = 10.7396941486331
$ws.Range("D6").Value = 4.261612520096079
$ws.Range("E6").Value = 13.42195239562019
$ws.Range("F6").Value = 21.302713342252
$ws.Range("G6").Value = 23.61558643264739
$ws.Range("H6").Value = 13.06798071634059
$ws.Range("L6").Value = 9.640358508111968
$ws.Range("M6").Value = 52.17819349555187
$ws.Range("O6").Value = 19.19152816663224

$ws.Range("C7").Value = 10.67019119023893
$ws.Range("D7").Value = 4.268246232625544
$ws.Range("E7").Value = 13.46899398040402
$ws.Range("F7").Value = 21.22181923842077
$ws.Range("G7").Value = 23.48186349582622
$ws.Range("H7").Value = 13.02141056947589
$ws.Range("L7").Value = 9.684691147216357
$ws.Range("M7").Value = 53.27356348455946
$ws.Range("O7").Value = 19.10466423588442

$ws.Range("C8").Value = 10.37878004377525
$ws.Range("D8").Value = 4.298661800000899
$ws.Range("E8").Value = 13.69252293202256
$ws.Range("F8").Value = 20.89545728689465
$ws.Range("G8").Value = 22.94368157619156
$ws.Range("H8").Value = 12.82634437854968
$ws.Range("L8").Value = 9.886175094579288
$ws.Range("M8").Value = 57.83781433370294
$ws.Range("O8").Value = 18.74404636696025

$ws.Range("C9").Value = 9.8640365922041
$ws.Range("D9").Value = 4.362124518167467
$ws.Range("E9").Value = 14.18576571462649
$ws.Range("F9").Value = 20.37082347298266
$ws.Range("G9").Value = 22.08897367757612
$ws.Range("H9").Value = 12.48178389355332
$ws.Range("L9").Value = 10.29932640819512
$ws.Range("M9").Value = 65.81128797368066
$ws.Range("O9").Value = 18.12107166992012

$ws.Range("C10").Value = 9.522397358691528
$ws.Range("D10").Value = 4.41082062132203
$ws.Range("E10").Value = 14.57917693383664
$ws.Range("F10").Value = 20.06085536484986
$ws.Range("G10").Value = 21.59726983427172
$ws.Range("H10").Value = 12.25222829778938
$ws.Range("L10").Value = 10.61121695556616
$ws.Range("M10").Value = 71.07097752607137
$ws.Range("O10").Value = 17.71736250994353

$ws.Range("C11").Value = 9.375436765107546
$ws.Range("D11").Value = 4.433368594600684
$ws.Range("E11").Value = 14.76432283398678
$ws.Range("F11").Value = 19.937392782542
$ws.Range("G11").Value = 21.4063651939672
$ws.Range("H11").Value = 12.15303022533389
$ws.Range("L11").Value = 10.75437480376046
$ws.Range("M11").Value = 73.33400213481909
$ws.Range("O11").Value = 17.54606870019764

$ws.Range("C12").Value = 9.321054205320618
$ws.Range("D12").Value = 4.441959253858561
$ws.Range("E12").Value = 14.83527265335666
$ws.Range("F12").Value = 19.8932627330675
$ws.Range("G12").Value = 21.33906390765895
$ws.Range("H12").Value = 12.11622871456059
$ws.Range("L12").Value = 10.80872543218088
$ws.Range("M12").Value = 74.17235701763013
$ws.Range("O12").Value = 17.483039511575

$ws.Range("C13").Value = 9.332709143371137
$ws.Range("D13").Value = 4.440106852793904
$ws.Range("E13").Value = 14.81995581307804
$ws.Range("F13").Value = 19.90264867571172
$ws.Range("G13").Value = 21.35333200995866
$ws.Range("H13").Value = 12.12412047677047
$ws.Range("L13").Value = 10.79701450619855
$ws.Range("M13").Value = 73.99262884599513
$ws.Range("O13").Value = 17.49653132535052

$ws.Range("C14").Value = 9.370936935259479
$ws.Range("D14").Value = 4.434074333224206
$ws.Range("E14").Value = 14.7701433867169
$ws.Range("F14").Value = 19.93370903694078
$ws.Range("G14").Value = 21.40072665130727
$ws.Range("H14").Value = 12.14998719291783
$ws.Range("L14").Value = 10.75884370305651
$ws.Range("M14").Value = 73.40334663916931
$ws.Range("O14").Value = 17.54084609856618

$ws.Range("C15").Value = 9.394519405696936
$ws.Range("D15").Value = 4.430385904471685
$ws.Range("E15").Value = 14.73973963448623
$ws.Range("F15").Value = 19.95307894144949
$ws.Range("G15").Value = 21.43041546395303
$ws.Range("H15").Value = 12.16593096036978
$ws.Range("L15").Value = 10.73547993382674
$ws.Range("M15").Value = 73.03997306636009
$ws.Range("O15").Value = 17.56823111959415

$ws.Range("C16").Value = 9.532176022028484
$ws.Range("D16").Value = 4.409354621593488
$ws.Range("E16").Value = 14.56719734291303
$ws.Range("F16").Value = 20.0692854326204
$ws.Range("G16").Value = 21.61042892143024
$ws.Range("H16").Value = 12.25881723419829
$ws.Range("L16").Value = 10.60188316248972
$ws.Range("M16").Value = 70.92047767776496
$ws.Range("O16").Value = 17.72881056556203

$ws.Range("C17").Value = 9.618824239864386
$ws.Range("D17").Value = 4.396550645699298
$ws.Range("E17").Value = 14.46289602987281
$ws.Range("F17").Value = 20.1451379080697
$ws.Range("G17").Value = 21.7294502739184
$ws.Range("H17").Value = 12.31714641667121
$ws.Range("L17").Value = 10.52022118598437
$ws.Range("M17").Value = 69.58705286841742
$ws.Range("O17").Value = 17.83052630843142

$ws.Range("C18").Value = 9.669454069882537
$ws.Range("D18").Value = 4.389223676799022
$ws.Range("E18").Value = 14.40348948582163
$ws.Range("F18").Value = 20.1904113909495
$ws.Range("G18").Value = 21.80097589710149
$ws.Range("H18").Value = 12.35118719426417
$ws.Range("L18").Value = 10.47337410086074
$ws.Range("M18").Value = 68.80789646954781
$ws.Range("O18").Value = 17.89018967161606

$ws.Range("C19").Value = 9.686731003063231
$ws.Range("D19").Value = 4.386749482389635
$ws.Range("E19").Value = 14.3834774069147
$ws.Range("F19").Value = 20.20601976754383
$ws.Range("G19").Value = 21.82571180430884
$ws.Range("H19").Value = 12.36279689340284
$ws.Range("L19").Value = 10.45753492966345
$ws.Range("M19").Value = 68.54199038493026
$ws.Range("O19").Value = 17.91058800755552

$ws.Range("C20").Value = 9.609518053183765
$ws.Range("D20").Value = 4.397909799796875
$ws.Range("E20").Value = 14.47393891715011
$ws.Range("F20").Value = 20.13689230085498
$ws.Range("G20").Value = 21.71646076983648
$ws.Range("H20").Value = 12.31088624421602
$ws.Range("L20").Value = 10.52890186541741
$ws.Range("M20").Value = 69.73026125878597
$ws.Range("O20").Value = 17.81957810896908

$ws.Range("C21").Value = 9.359673631782698
$ws.Range("D21").Value = 4.435844850011139
$ws.Range("E21").Value = 14.78475212743938
$ws.Range("F21").Value = 19.92451387335156
$ws.Range("G21").Value = 21.38666801717476
$ws.Range("H21").Value = 12.1423687222238
$ws.Range("L21").Value = 10.7700519359228
$ws.Range("M21").Value = 73.57693758674036
$ws.Range("O21").Value = 17.52777945115291

$ws.Range("C22").Value = 9.203807678618242
$ws.Range("D22").Value = 4.460939910398037
$ws.Range("E22").Value = 14.99275448598344
$ws.Range("F22").Value = 19.80104645619437
$ws.Range("G22").Value = 21.20034943867814
$ws.Range("H22").Value = 12.03668376305938
$ws.Range("L22").Value = 10.92845391600188
$ws.Range("M22").Value = 75.98256051342474
$ws.Range("O22").Value = 17.34780474210815

$ws.Range("C23").Value = 9.28629789745942
$ws.Range("D23").Value = 4.447520077266306
$ws.Range("E23").Value = 14.88131068040712
$ws.Range("F23").Value = 19.86550668573786
$ws.Range("G23").Value = 21.2970240537099
$ws.Range("H23").Value = 12.09267872097469
$ws.Range("L23").Value = 10.84385282307847
$ws.Range("M23").Value = 74.7085349918659
$ws.Range("O23").Value = 17.44285787275052

$ws.Range("C24").Value = 9.61372284673692
$ws.Range("D24").Value = 4.397295219324128
$ws.Range("E24").Value = 14.46894468701533
$ws.Range("F24").Value = 20.14061496166634
$ws.Range("G24").Value = 21.7223236939141
$ws.Range("H24").Value = 12.3137148911941
$ws.Range("L24").Value = 10.52497701083075
$ws.Range("M24").Value = 69.66555584703099
$ws.Range("O24").Value = 17.82452410142037

$ws.Range("C25").Value = 9.997030542525803
$ws.Range("D25").Value = 4.344573972859826
$ws.Range("E25").Value = 14.04669488391836
$ws.Range("F25").Value = 20.49987887521158
$ws.Range("G25").Value = 22.29731586684986
$ws.Range("H25").Value = 12.57088828868422
$ws.Range("L25").Value = 10.18590306111382
$ws.Range("M25").Value = 63.75929927669386
$ws.Range("O25").Value = 18.28032617716875
